# Updates cryptos list (prices/volumes, plus a couple of row reorderings)
# as captured by the latest GitHub Actions scrape.
#
# Note: several "Price" values in column D (e.g. "628.05") look like plain
# numbers to Excel's auto-detection, which would otherwise silently convert
# them to floating point numbers (losing exact text like trailing zeros,
# e.g. "1.00" -> 1). To preserve them as text (matching the source data,
# which always stores these as inline strings) we temporarily force the
# cell to Text format before assigning the value, then clear the
# formatting again so the cell keeps its original (default) style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.065.60'
$ws.Range('E2').Value = '  +2.15%  '
$ws.Range('D3').Value = '3.817.65'
$ws.Range('E3').Value = '  +0.89%  '
$ws.Range('E4').Value = '  +0.35%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '628.05'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +5.10%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '165.36'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.47%  '
$ws.Range('D7').Value = '3.815.60'
$ws.Range('E7').Value = '  +0.99%  '
$ws.Range('E8').Value = '  -0.19%  '
$ws.Range('E9').Value = '  +0.99%  '
$ws.Range('E10').Value = '  +2.35%  '
$ws.Range('E11').Value = '  +1.14%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.61'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +3.07%  '
$ws.Range('E13').Value = '  +1.19%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.03'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.24%  '
$ws.Range('D15').Value = '4.454.02'
$ws.Range('E15').Value = '  +0.73%  '
$ws.Range('D16').Value = '3.815.38'
$ws.Range('E16').Value = '  +0.56%  '
$ws.Range('D17').Value = '69.002.61'
$ws.Range('E17').Value = '  +1.98%  '
$ws.Range('E18').Value = '  -1.35%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.12'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.47%  '
$ws.Range('E20').Value = '  -0.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '466.37'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +1.00%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.69'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.65%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.709'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.95%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000152'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +4.99%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '83.74'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +1.36%  '
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.96'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.15%  '
$ws.Range('B27').Value = 'Fetch.AI'
$ws.Range('C27').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.15'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +2.99%  '
$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.03'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.50%  '
$ws.Range('B29').Value = 'Dai'
$ws.Range('C29').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.02%  '
$ws.Range('D30').Value = '3.965.81'
$ws.Range('E30').Value = '  +0.85%  '
$ws.Range('E31').Value = '  +1.63%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.22'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.55%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.29'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.46%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '29.21'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.87%  '
$ws.Range('E35').Value = '  +0.08%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '9.09'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +1.33%  '
$ws.Range('E37').Value = '  +2.92%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.148'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +7.29%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.41'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +6.44%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.91'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +3.21%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.980'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.48%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.999'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.01%  '
$ws.Range('E43').Value = '  +0.01%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '156.95'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +3.47%  '
$ws.Range('E45').Value = '  +6.06%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.299'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.84%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '46.79'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.49%  '
$ws.Range('E48').Value = '  -3.13%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.45'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +1.68%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.90'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +2.79%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.000281'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +14.25%  '
